{"js": "// Update the \"id\" row of the \"resource\" table (GEO-RES schema doc):\n//  - rename tag \"id\" -> \"resourceId\"\n//  - update field/description/example text\n//  - add a REGEX format hint and a two-line example\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The second table (index 1) is the \"resource\" table; its second row\n// (index 1, just below the header row) describes the \"id\" field.\nconst table = tables.items[1];\nconst row = 1;\n\nconst tagCell = table.getCell(row, 0);\nconst fieldCell = table.getCell(row, 1);\nconst formatCell = table.getCell(row, 2);\nconst descCell = table.getCell(row, 4);\nconst exampleCell = table.getCell(row, 5);\n\n// \\u000b (vertical tab) maps to a <w:br/> line break when inserted via\n// Office.js, splitting the text into separate runs around the break.\ntagCell.body.insertText(\"resourceId\", Word.InsertLocation.replace);\n\nfieldCell.body.insertText(\n  \"Identifiant de la ressource partag\u00e9\",\n  Word.InsertLocation.replace\n);\n\nformatCell.body.insertText(\n  \"string\\u000b(REGEX: ^([\\\\w-]+\\\\.){3,4}resource(\\\\.[\\\\w-]+){1,2}$)\",\n  Word.InsertLocation.replace\n);\n\ndescCell.body.insertText(\n  \"A valoriser avec l'identifiant partag\u00e9 unique de la ressource engag\u00e9e, norm\u00e9 comme suit :\\u000b\" +\n    \"{orgID}.resource.{ID unique de la ressource partag\u00e9e}\\u000b\" +\n    \"OU - uniquement dans le cas o\u00f9 un ID unique de ressource ne peut pas \u00eatre garanti par l'organisation propri\u00e9taire :\\u000b\" +\n    \"{orgID}.resource.{sendercaseId}.{n\u00b0 d\\u2019ordre chronologique de la ressource}\",\n  Word.InsertLocation.replace\n);\n\nexampleCell.body.insertText(\n  \"fr.health.samu770.resource.VLM250\\u000bfr.health.samu440.resource.DRFR15DDXAAJJJ0000.1\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Update the \"id\" row of the \"resource\" table (GEO-RES schema doc):\n#  - rename tag \"id\" -> \"resourceId\"\n#  - update field/description/example text\n#  - add a REGEX format hint and a two-line example\n$d = $word.ActiveDocument\n\n# The second table is the \"resource\" table; its second row (row 2,\n# just below the header row) describes the \"id\" field. Word COM\n# table/cell indices are 1-based.\n$table = $d.Tables.Item(2)\n$row = 2\n\n# [char]11 is a manual line break (Shift+Enter) -> serializes as <w:br/>\n# splitting the text into separate <w:t> runs around the break, same as\n# typing Shift+Enter within a single run in Word.\n$nl = [char]11\n\n$table.Cell($row, 1).Range.Text = \"resourceId\"\n\n$table.Cell($row, 2).Range.Text = \"Identifiant de la ressource partag\u00e9\"\n\n$table.Cell($row, 3).Range.Text = \"string\" + $nl + \"(REGEX: ^([\\w-]+\\.){3,4}resource(\\.[\\w-]+){1,2}$)\"\n\n$descLine1 = \"A valoriser avec l'identifiant partag\u00e9 unique de la ressource engag\u00e9e, norm\u00e9 comme suit :\"\n$descLine2 = \"{orgID}.resource.{ID unique de la ressource partag\u00e9e}\"\n$descLine3 = \"OU - uniquement dans le cas o\u00f9 un ID unique de ressource ne peut pas \u00eatre garanti par l'organisation propri\u00e9taire :\"\n$descLine4 = \"{orgID}.resource.{sendercaseId}.{n\u00b0 d\u2019ordre chronologique de la ressource}\"\n$table.Cell($row, 5).Range.Text = $descLine1 + $nl + $descLine2 + $nl + $descLine3 + $nl + $descLine4\n\n$table.Cell($row, 6).Range.Text = \"fr.health.samu770.resource.VLM250\" + $nl + \"fr.health.samu440.resource.DRFR15DDXAAJJJ0000.1\"\n"}
